# Added Assignment 3 to WBA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grab formatting references from existing rows so new rows inherit the
# --- exact same look (fill colors come from the theme, so copy them from
# --- cells that already carry the fill we need). ---

# B4:B8 use the "task name" fill (s=2 -> fillId 4, no alignment override)
$taskFill = $ws.Range("B4").Interior.Color
# E3:F3 use fillId 2 + centered + wrap (same definition reused for the new
# section header E9:F9)
$sectionFill = $ws.Range("E3").Interior.Color

# ---------------------------------------------------------------------
# Row 9: new section header "Code and Documentation Have Ready By"
# ---------------------------------------------------------------------
$ws.Range("E9:F9").Merge()
$ws.Range("E9").Value = "Code and Documentation Have Ready By"
$ws.Range("E9:F9").Interior.Color = $sectionFill
$ws.Range("E9:F9").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E9:F9").WrapText = $true
$ws.Rows("9:9").RowHeight = 31.5

# ---------------------------------------------------------------------
# Rows 10-13: four new Assignment 3 tasks
# ---------------------------------------------------------------------
$tasks = @(
    @{ Row = 10; Task = "Going to Town";    Difficulty = "Hard";   Person = "Rahn" },
    @{ Row = 11; Task = "Mambo Marie";      Difficulty = "Medium"; Person = "Rahn" },
    @{ Row = 12; Task = "Ending the Game";  Difficulty = "Medium"; Person = "Arie" },
    @{ Row = 13; Task = "Guns";             Difficulty = "Medium"; Person = "Arie" }
)

foreach ($t in $tasks) {
    $r = $t.Row
    $ws.Cells.Item($r, 2).Value = $t.Task
    $ws.Cells.Item($r, 2).Interior.Color = $taskFill

    $ws.Cells.Item($r, 3).Value = $t.Difficulty
    $ws.Cells.Item($r, 4).Value = $t.Person

    $ws.Range("E" + $r + ":F" + $r).Merge()
    $ws.Range("E" + $r).Value = "8th of June"
    $ws.Range("E" + $r + ":F" + $r).HorizontalAlignment = -4108   # xlCenter
}

# ---------------------------------------------------------------------
# Row 14: Recommendations Report / Both (Half/Half)
# ---------------------------------------------------------------------
$ws.Range("B14").Value = "Recommendations Report"
$ws.Range("B14").Interior.Color = $taskFill
$ws.Range("B14").WrapText = $true

$ws.Range("D14").Value = "Both (Half/Half)"
$ws.Range("D14").WrapText = $true

$ws.Range("E14:F14").Merge()
$ws.Range("E14").Value = "8th of June"
$ws.Range("E14:F14").HorizontalAlignment = -4108   # xlCenter
$ws.Rows("14:14").RowHeight = 45

# ---------------------------------------------------------------------
# Column width tweaks (best effort to match final layout)
# ---------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 10.71
$ws.Columns("E:E").ColumnWidth = 15.14
$ws.Columns("F:F").ColumnWidth = 12.43

# Cursor position, like the source file
$ws.Range("P11").Select()
